$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update 想去人数 (F column) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 250
$wsExhibit.Range("F3").Value = 78
$wsExhibit.Range("F4").Value = 823
$wsExhibit.Range("F5").Value = 520

# Sheet "全部类型" (sheet4): update 想去人数 (F column) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 250
$wsAll.Range("F3").Value = 78
$wsAll.Range("F4").Value = 823
$wsAll.Range("F6").Value = 520
